$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RESTAURANT table: rename total_tables -> restaurant_location (E6)
$ws.Range("E6").Value = "restaurant_location"

# New STAFF RESTAURANTS table in column O (rows 3-6)
$ws.Range("O3").Value = "STAFF RESTAURANTS"
$ws.Range("O4").Value = "restaurant_id"
$ws.Range("O5").Value = "user_id (restaurant users)"
$ws.Range("O6").Value = "account_type"

# Widen the new column to fit its header/content
$ws.Columns.Item(15).ColumnWidth = 24.26

# Update selection to mirror the saved cursor position
$ws.Range("I20").Select() | Out-Null
